# Rotate the species-record data for rows 2-5 (the fields that vary per
# record: Id, Taxonsorteringsordning, Rodlistade, TaxonId, Artnamn,
# Vetenskapligt namn, Auktor, Ost, Nord) so that:
#   new row 2 <- old row 5
#   new row 3 <- old row 2
#   new row 4 <- old row 3
#   new row 5 <- old row 4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

# Snapshot the current values of rows 2-5 for the columns that change.
# (Value2 is used since it round-trips numbers/strings faithfully here.)
$orig = @{}
foreach ($col in $cols) {
    $orig[$col] = @{}
    for ($r = 2; $r -le 5; $r++) {
        $orig[$col][$r] = $ws.Range("$col$r").Value2
    }
}

# Source row mapping: destination row -> source row
$srcRow = @{ 2 = 5; 3 = 2; 4 = 3; 5 = 4 }

foreach ($col in $cols) {
    foreach ($destRow in 2..5) {
        $s = $srcRow[$destRow]
        $ws.Range("$col$destRow").Value2 = $orig[$col][$s]
    }
}
